$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at R:T, shifting the existing R:AE columns
# (general_college_subjects.arts onward) to the right.
$ws.Range("R1:T2").Insert(-4161)

# Set the headers for the newly inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# Copy the style of the existing header cells onto the new header cells.
$ws.Range("U1").Copy() | Out-Null
$ws.Range("R1:T1").PasteSpecial(-4122) | Out-Null

# Set the data values for the newly inserted columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 0

# Lowercase the "Unknown" values in row 2.
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"
